$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect so the cell values can be written,
# then restore protection once all edits are applied.
$ws.Unprotect()

# Refresh the "as of" date in the confidentiality footer (A42):
# 2021-04-22 -> 2021-04-23
$ws.Range("A42").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# Updated Weight (column D) and Percent Change (column E) figures for each
# holding row, row number -> (Weight, PercentChange). A blank Weight means
# that column is unchanged for that row (only the Total's Percent Change
# moved on row 39).
$updates = @(
    @(2,  "0.06254140732772734",  "0.01803850234955262"),
    @(3,  "0.05669868870089519",  "0.01547614418478038"),
    @(4,  "0.2882887014243129",   "0.02596359743040688"),
    @(5,  "0.03647747187829261",  "0.009622126054686708"),
    @(6,  "0.03261227352327552",  "0.01067803313266258"),
    @(7,  "0.02924181350260851",  "0.01913550926240082"),
    @(8,  "0.02913400287627287",  "0.002058360576340901"),
    @(9,  "0.02463461788188056",  "0.001646738741319043"),
    @(10, "0.02483083763124401",  "0.02104753786869806"),
    @(11, "0.02288098654880547",  "0.01554701200593556"),
    @(12, "0.02198895833990155",  "0.02137643378519294"),
    @(13, "0.02216136715338718",  "0.001422630772597966"),
    @(14, "0.02226157150853964",  "-0.005125157840005889"),
    @(15, "0.02126184290911432",  "0.007559260872270324"),
    @(16, "0.0218564548911741",   "0.009507237605285601"),
    @(17, "0.02113000087527237",  "0.009651502504173459"),
    @(18, "0.01655389175416267",  "-0.05322039315966121"),
    @(19, "0.01679111922642679",  "-0.0001838235294117085"),
    @(20, "0.01578576860047477",  "0.0003491620111730764"),
    @(21, "0.01584110697922782",  "0.005427899402931091"),
    @(22, "0.01586712263139062",  "0.01349192013227918"),
    @(23, "0.01560321809207821",  "0.0005510653930933529"),
    @(24, "0.01486507907986607",  "0.001275510204081565"),
    @(25, "0.0141928610845749",   "0.007961165048543606"),
    @(26, "0.01497553536574369",  "-0.001538461538461489"),
    @(27, "0.01281623623623182",  "0.007818547762811434"),
    @(28, "0.01335396653813906",  "0.005943536404160232"),
    @(29, "0.01448234033279311",  "0.0005175983436851439"),
    @(30, "0.01320536866053124",  "0.001335648457325922"),
    @(31, "0.0128817163098958",   "-0.001643047853768564"),
    @(32, "0.01334459649392788",  "0.01208541572012711"),
    @(33, "0.01289097611829273",  "-0.001026167265264277"),
    @(34, "0.006548117602212301", "0.02794565748051392"),
    @(35, "0.005608577757367004", "-0.006348519988993218"),
    @(36, "0.005782309400623626", "0.01433637091546869"),
    @(37, "0.005592152621043885", "0.01685426481893981"),
    @(38, "0.005016942142291948", "0.02421392630353103"),
    @(39, $null,                  "0.01300484971441218")
)

foreach ($row in $updates) {
    $rowNum = $row[0]
    $weight = $row[1]
    $pctChange = $row[2]

    if ($weight -ne $null) {
        $ws.Range("D$rowNum").Value = [double]$weight
    }
    $ws.Range("E$rowNum").Value = [double]$pctChange
}

# Restore sheet protection. The original legacy password hash cannot be
# reversed into its plaintext, so re-lock the sheet without a password.
$ws.Protect()
